# Merge the multiple single-word/space runs in the Title, Author and
# Abstract paragraphs into one run each, per the target diff. Using
# Find/Execute over each paragraph's own range (with MatchCase so the
# search is unambiguous) replaces the whole matched span - which covers
# every run in the paragraph - with a single new run containing the
# full paragraph text.

$d = $word.ActiveDocument

# "Questions:" " " "Hypothesis" " " "Testing"  ->  "Questions: Hypothesis Testing"
$title = $d.Paragraphs(1).Range
$title.Find.Execute("Questions: Hypothesis Testing", $true, $false, $false, $false, $false, $true, 1, $false, "Questions: Hypothesis Testing", 2)

# "Ellie" " " "Trace"  ->  "Ellie Trace"
$author = $d.Paragraphs(2).Range
$author.Find.Execute("Ellie Trace", $true, $false, $false, $false, $false, $true, 1, $false, "Ellie Trace", 2)

# "A" " " "selection" " " ... " " "Testing."  ->  one run with the full sentence
$abstract = $d.Paragraphs(4).Range
$abstract.Find.Execute("A selection of questions for the study guide on Hypothesis Testing.", $true, $false, $false, $false, $false, $true, 1, $false, "A selection of questions for the study guide on Hypothesis Testing.", 2)
